$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6571.7144
$ws.Range("J51").Value = 4400.4
$ws.Range("L51").Value = 4400.4
$ws.Range("N51").Value = -5368.4

$ws.Range("H86").Value = 34214
$ws.Range("J86").Value = 34214
$ws.Range("L86").Value = 34214
$ws.Range("N86").Value = -36460

$ws.Range("H89").Value = 34214
$ws.Range("J89").Value = 34214
$ws.Range("L89").Value = 171070
$ws.Range("N89").Value = -182302

$ws.Range("H92").Value = 550
$ws.Range("I92").Value = 550
$ws.Range("K92").Value = 550
$ws.Range("M92").Value = 698

$ws.Range("H103").Value = 434.66666
$ws.Range("I103").Value = 452
$ws.Range("J103").Value = 400
$ws.Range("K103").Value = 1356
$ws.Range("L103").Value = 1200
$ws.Range("M103").Value = -770
$ws.Range("N103").Value = -2372

$ws.Range("H129").Value = 257379.97
$ws.Range("J129").Value = 286757.7
$ws.Range("L129").Value = 860273.1000000001
$ws.Range("N129").Value = -870273.1000000001

$ws.Range("H135").Value = 20005992
$ws.Range("I135").Value = 1047.6842
$ws.Range("K135").Value = 9429.157799999999
$ws.Range("M135").Value = -6894.157799999999

$ws.Range("H137").Value = 95867.92999999999
$ws.Range("I137").Value = 119638.85
$ws.Range("J137").Value = 6066.6665
$ws.Range("K137").Value = 358916.55
$ws.Range("L137").Value = 18199.9995
$ws.Range("M137").Value = -356366.55
$ws.Range("N137").Value = -23299.9995

$ws.Range("H138").Value = 3805.3647
$ws.Range("J138").Value = 4085.7937
$ws.Range("L138").Value = 12257.3811
$ws.Range("N138").Value = -22537.3811

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2074
$ws.Range("I2").Value = 1987.5555
$ws.Range("K2").Value = 1987.5555
$ws.Range("M2").Value = -1874.5555

$ws.Range("H32").Value = 16459.947
$ws.Range("I32").Value = 11243.372
$ws.Range("K32").Value = 11243.372
$ws.Range("M32").Value = -10956.372

$ws.Range("H45").Value = 4500.7144
$ws.Range("I45").Value = 4556.2144
$ws.Range("J45").Value = 4389.7144
$ws.Range("K45").Value = 4556.2144
$ws.Range("L45").Value = 4389.7144
$ws.Range("M45").Value = -4179.2144
$ws.Range("N45").Value = -5143.7144

$ws.Range("H61").Value = 723485.1
$ws.Range("I61").Value = 1003307.1
$ws.Range("K61").Value = 1003307.1
$ws.Range("M61").Value = -1003095.1

$ws.Range("H74").Value = 28574284
$ws.Range("I74").Value = 43480732
$ws.Range("K74").Value = 43480732
$ws.Range("M74").Value = -43479858

$ws.Range("H77").Value = 28574284
$ws.Range("I77").Value = 43480732
$ws.Range("K77").Value = 217403660
$ws.Range("M77").Value = -217399292

$ws.Range("H88").Value = 112740.89
$ws.Range("J88").Value = 168228
$ws.Range("L88").Value = 168228
$ws.Range("N88").Value = -169040

$ws.Range("H91").Value = 112740.89
$ws.Range("J91").Value = 168228
$ws.Range("L91").Value = 168228
$ws.Range("N91").Value = -171036

$ws.Range("H97").Value = 270.3889
$ws.Range("I97").Value = 236.83333
$ws.Range("J97").Value = 337.5
$ws.Range("K97").Value = 236.83333
$ws.Range("L97").Value = 337.5
$ws.Range("M97").Value = 259.16667
$ws.Range("N97").Value = -1329.5

$ws.Range("H110").Value = 1717.9615
$ws.Range("I110").Value = 1080.55
$ws.Range("J110").Value = 3842.6667
$ws.Range("K110").Value = 1080.55
$ws.Range("L110").Value = 3842.6667
$ws.Range("M110").Value = 964.45
$ws.Range("N110").Value = -7932.6667

$ws.Range("H116").Value = 2074
$ws.Range("I116").Value = 1987.5555
$ws.Range("K116").Value = 1987.5555
$ws.Range("M116").Value = 306.4445000000001

$ws.Range("H122").Value = 1580.8182
$ws.Range("I122").Value = 1688.1177
$ws.Range("K122").Value = 5064.3531
$ws.Range("M122").Value = -2614.3531

$ws.Range("H136").Value = 723485.1
$ws.Range("I136").Value = 1003307.1
$ws.Range("K136").Value = 3009921.3
$ws.Range("M136").Value = -3007371.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2074
$ws.Range("I3").Value = 1987.5555
$ws.Range("K3").Value = 1987.5555
$ws.Range("M3").Value = -1873.5555

$ws.Range("H94").Value = 1042.8966
$ws.Range("I94").Value = 474.1111
$ws.Range("J94").Value = 1973.6364
$ws.Range("K94").Value = 474.1111
$ws.Range("L94").Value = 1973.6364
$ws.Range("M94").Value = -23.11110000000002
$ws.Range("N94").Value = -2875.6364

$ws.Range("H105").Value = 4611926
$ws.Range("I105").Value = 7577250.5
$ws.Range("K105").Value = 7577250.5
$ws.Range("M105").Value = -7575503.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1085.5
$ws.Range("I16").Value = 1262.8182
$ws.Range("K16").Value = 1262.8182
$ws.Range("M16").Value = -975.8181999999999

$ws.Range("H31").Value = 7040.7173
$ws.Range("I31").Value = 3914.75
$ws.Range("J31").Value = 8144
$ws.Range("K31").Value = 3914.75
$ws.Range("L31").Value = 8144
$ws.Range("M31").Value = -3619.75
$ws.Range("N31").Value = -8734

$ws.Range("H34").Value = 7040.7173
$ws.Range("I34").Value = 3914.75
$ws.Range("J34").Value = 8144
$ws.Range("K34").Value = 3914.75
$ws.Range("L34").Value = 8144
$ws.Range("M34").Value = -3712.75
$ws.Range("N34").Value = -8548

$ws.Range("H86").Value = 7513
$ws.Range("I86").Value = 1900
$ws.Range("K86").Value = 1900
$ws.Range("M86").Value = -777

$ws.Range("H89").Value = 7513
$ws.Range("I89").Value = 1900
$ws.Range("K89").Value = 9500
$ws.Range("M89").Value = -3884

$ws.Range("H107").Value = 2398.2666
$ws.Range("I107").Value = 1740.2
$ws.Range("K107").Value = 1740.2
$ws.Range("M107").Value = 179.8

$ws.Range("H113").Value = 1085.5
$ws.Range("I113").Value = 1262.8182
$ws.Range("K113").Value = 1262.8182
$ws.Range("M113").Value = 907.1818000000001

$ws.Range("H134").Value = 1015.3333
$ws.Range("I134").Value = 964.3158
$ws.Range("K134").Value = 2892.9474
$ws.Range("M134").Value = -357.9474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1374.8125
$ws.Range("I122").Value = 559.6
$ws.Range("J122").Value = 1525.7778
$ws.Range("K122").Value = 5036.400000000001
$ws.Range("L122").Value = 13732.0002
$ws.Range("M122").Value = -2586.400000000001
$ws.Range("N122").Value = -18632.0002

$ws.Range("H131").Value = 769.29
$ws.Range("J131").Value = 769.29
$ws.Range("L131").Value = 2307.87
$ws.Range("N131").Value = -12387.87

$ws.Range("H137").Value = 13339768
$ws.Range("I137").Value = 2771.1428
$ws.Range("J137").Value = 18526378
$ws.Range("K137").Value = 8313.428400000001
$ws.Range("L137").Value = 55579134
$ws.Range("M137").Value = -3213.428400000001
$ws.Range("N137").Value = -55589334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3476.25
$ws.Range("I82").Value = 2300.6667
$ws.Range("K82").Value = 2300.6667
$ws.Range("M82").Value = -1939.6667

$ws.Range("H85").Value = 3476.25
$ws.Range("I85").Value = 2300.6667
$ws.Range("K85").Value = 2300.6667
$ws.Range("M85").Value = -1052.6667

$ws.Range("H93").Value = 2078.7144
$ws.Range("I93").Value = 2100.182
$ws.Range("K93").Value = 2100.182
$ws.Range("M93").Value = -852.1819999999998

$ws.Range("H136").Value = 2917.111
$ws.Range("I136").Value = 2420.8
$ws.Range("K136").Value = 7262.400000000001
$ws.Range("M136").Value = -4712.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4526.5557
$ws.Range("I62").Value = 4163
$ws.Range("K62").Value = 4163
$ws.Range("M62").Value = -3539

$ws.Range("H65").Value = 4526.5557
$ws.Range("I65").Value = 4163
$ws.Range("K65").Value = 20815
$ws.Range("M65").Value = -17695

$ws.Range("H136").Value = 21741742
$ws.Range("I136").Value = 29413058
$ws.Range("K136").Value = 88239174
$ws.Range("M136").Value = -88236624
